$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B ("ASIN" and everything to its right shifts over by one).
$ws.Columns("B:B").Insert()

# New column header + data.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "Week_Start_Date"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
}

# Week labels lose their leading zero (W01 -> W1 ... W09 -> W9; W10-W16 already match).
for ($w = 1; $w -le 16; $w++) {
    $row = $w + 1
    $ws.Range("A$row").Value = "W$w"
}

# is_holiday_week (now column J) becomes a boolean column instead of numeric 0/1.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}
